$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New product row (id 11, "Taladro Inalambrico") and a run of blank id rows
# down to 20, replicating the order the values were originally typed in so
# that shared-string allocation order matches.

$ws.Range("A12").Value = "11"
$ws.Range("F12").Value = "edward quevedo"
$ws.Range("C12").Value = "Taladro Inalambrico"

$ws.Range("A13").Value = "12"

$ws.Range("A16").Value = "15"
$ws.Range("A17").Value = "16"
$ws.Range("A18").Value = "17"
$ws.Range("A19").Value = "18"
$ws.Range("A20").Value = "19"

$ws.Range("A14").Value = "13"
$ws.Range("A15").Value = "14"

$ws.Range("A21").Value = "20"

$ws.Range("D12").Value = "130.000"
$ws.Range("E12").Value = "Taladro Inalámbrico 21 Voltios, Percutor Con 2 Baterías, Copas, Puntas Y Extensión"
$ws.Range("G12").Value = "3224336215"

$imgUrl = "https://lh3.googleusercontent.com/pw/AP1GczObtmbU_tsVak1pMSlQuN4-4Onz8TKrQnKBd-dycgx41VHwpHSVWG2njJT7U9vEMJejJa1fdxxmZc7CAuM4pVHM37jFFDfRJchHoHvxUA3iQGsdSAIgmfdNFRoLJsJjMa-xkMW_15cmerYDk3qJjnt-5g=w475-h633-s-no-gm?authuser=0"
$ws.Range("B12").Value = $imgUrl
$ws.Range("I12").Value = $imgUrl

# oferta column re-uses the existing "no" shared string
$ws.Range("H12").Value = "no"

# Wrap the long description text and give the row extra height for it
$ws.Range("E12").WrapText = $true
$ws.Rows.Item(12).RowHeight = 30

# The id column cells carry the sheet's text number-format (same style as the
# rest of column A); stamp it explicitly on the price column's blank
# continuation cells (D13:D20) so they persist even without a value.
$ws.Range("D13").NumberFormat = $ws.Range("A13").NumberFormat
$ws.Range("D14").NumberFormat = $ws.Range("A14").NumberFormat
$ws.Range("D15").NumberFormat = $ws.Range("A15").NumberFormat
$ws.Range("D16").NumberFormat = $ws.Range("A16").NumberFormat
$ws.Range("D17").NumberFormat = $ws.Range("A17").NumberFormat
$ws.Range("D18").NumberFormat = $ws.Range("A18").NumberFormat
$ws.Range("D19").NumberFormat = $ws.Range("A19").NumberFormat
$ws.Range("D20").NumberFormat = $ws.Range("A20").NumberFormat

# Update the view: scrolled over to column G, active cell on J16
$ws.Range("J16").Select()
